$d = $word.ActiveDocument

# Locate the paragraph that holds "Ver no Jupiter Salvar em pdf Salvar em docx"
# and the paragraph that holds the copyright notice ("... Powered by Jekyll and
# Github pages ..."). Those two paragraphs, plus the blank paragraph that
# immediately precedes the first one (right after the "LOB1018: Física I
# (Requisito fraco)" requirements line), are removed entirely - this is the
# footer block that the Jekyll site build no longer emits for this page.

$count = $d.Paragraphs.Count
$startPara = $null
$endPara = $null

for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $startPara = $d.Paragraphs.Item($i - 1)
    }
    if ($t -like "*Powered by Jekyll and Github pages*") {
        $endPara = $p
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $deleteRange.Delete()
}
